# This script updates the cryptocurrency price/volume data in the active worksheet
# to match the latest scrape, including two pairs of rows whose rank order swapped
# (NEARProtocol/OKB at rows 38-39, and Stacks/Monero at rows 49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") sometimes holds purely numeric-looking text (e.g. "0.999").
# Assigning such a string straight to Range.Value lets Excel auto-convert it to a
# real number, which would not match the original text-cell formatting. Prefixing
# with a single quote forces Excel to keep/store it as text, exactly as before.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
}

$ws.Range("D2").Value = "66.950.78"
$ws.Range("E2").Value = "  +7.52%  "

$ws.Range("D3").Value = "3.866.68"
$ws.Range("E3").Value = "  +10.71%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.76%  "

Set-TextValue $ws.Range("D5") "424.62"
$ws.Range("E5").Value = "  +7.81%  "

Set-TextValue $ws.Range("D6") "131.39"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("D7").Value = "3.860.82"
$ws.Range("E7").Value = "  +6.06%  "

$ws.Range("E8").Value = "  +2.79%  "

Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("E10").Value = "  +5.21%  "

$ws.Range("E11").Value = "  +5.19%  "

$ws.Range("E12").Value = "  +16.65%  "

$ws.Range("E13").Value = "  +2.27%  "

Set-TextValue $ws.Range("D14") "10.29"
$ws.Range("E14").Value = "  +9.11%  "

$ws.Range("D15").Value = "4.470.47"
$ws.Range("E15").Value = "  +8.59%  "

Set-TextValue $ws.Range("D16") "15.91"
$ws.Range("E16").Value = "  +24.74%  "

$ws.Range("D17").Value = "3.867.01"
$ws.Range("E17").Value = "  +11.04%  "

$ws.Range("E18").Value = "  +0.12%  "

Set-TextValue $ws.Range("D19") "20.05"
$ws.Range("E19").Value = "  +4.76%  "

$ws.Range("D20").Value = "67.169.11"
$ws.Range("E20").Value = "  +7.56%  "

$ws.Range("E21").Value = "  +5.20%  "

Set-TextValue $ws.Range("D22") "414.64"
$ws.Range("E22").Value = "  +2.55%  "

Set-TextValue $ws.Range("D23") "14.99"
$ws.Range("E23").Value = "  +4.40%  "

Set-TextValue $ws.Range("D24") "84.45"
$ws.Range("E24").Value = "  +2.78%  "

Set-TextValue $ws.Range("D25") "3.04"
$ws.Range("E25").Value = "  +6.06%  "

Set-TextValue $ws.Range("D26") "37.79"
$ws.Range("E26").Value = "  +11.08%  "

Set-TextValue $ws.Range("D27") "9.94"
$ws.Range("E27").Value = "  +10.71%  "

Set-TextValue $ws.Range("D28") "3.25"
$ws.Range("E28").Value = "  +5.60%  "

Set-TextValue $ws.Range("D29") "5.29"
$ws.Range("E29").Value = "  +3.44%  "

Set-TextValue $ws.Range("D30") "9.14"
$ws.Range("E30").Value = "  +38.14%  "

Set-TextValue $ws.Range("D31") "721.99"
$ws.Range("E31").Value = "  +10.67%  "

Set-TextValue $ws.Range("D32") "13.23"
$ws.Range("E32").Value = "  +9.26%  "

Set-TextValue $ws.Range("D33") "0.122"
$ws.Range("E33").Value = "  +10.07%  "

$ws.Range("E34").Value = "  +4.96%  "

Set-TextValue $ws.Range("D35") "0.998"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("E36").Value = "  -2.06%  "

Set-TextValue $ws.Range("D37") "38.89"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D38") "55.49"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D39") "5.45"
$ws.Range("E39").Value = "  +31.34%  "

$ws.Range("D40").Value = "0.0₃0754"
$ws.Range("E40").Value = "  +26.52%  "

Set-TextValue $ws.Range("D41") "0.0463"
$ws.Range("E41").Value = "  +4.46%  "

Set-TextValue $ws.Range("D42") "2.89"
$ws.Range("E42").Value = "  +4.80%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("E44").Value = "  +6.77%  "

$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("E46").Value = "  +3.95%  "

Set-TextValue $ws.Range("D47") "0.313"
$ws.Range("E47").Value = "  +13.48%  "

$ws.Range("E48").Value = "  +4.27%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D49") "141.44"
$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D50") "2.82"
$ws.Range("E50").Value = "  +4.35%  "

$ws.Range("E51").Value = "  +2.49%  "
